$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "64.246.22"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -3.63%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.156.59"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -2.98%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "606.08"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.48"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -6.59%  "

$ws.Range("E7").Value = "  -0.04%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.150.49"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -3.25%  "

$ws.Range("E9").Value = "  -3.89%  "

$ws.Range("E10").Value = "  -6.54%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.52"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -6.61%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.476"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.97%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000251"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -7.42%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.77"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -9.14%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.667.45"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.34%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.205.46"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.82%  "

$ws.Range("E17").Value = "  +0.61%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.147.73"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.43%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.75%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "481.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -5.29%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.80"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.73%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.710"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -5.65%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.73"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.96%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "13.72"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -7.47%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "83.74"
$c.Style = "Normal"

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -5.09%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.47"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.83%  "

$ws.Range("E29").Value = "  -8.10%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.80"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -30.42%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.75"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.57%  "

$ws.Range("E33").Value = "  +0.03%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "26.25"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -6.74%  "

$ws.Range("E35").Value = "  -4.67%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "54.26"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -6.62%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0₃0719"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -10.65%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "452.89"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -8.56%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -13.03%  "

$ws.Range("E41").Value = "  -7.81%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "8.45"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.95%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -8.38%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.846.85"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.23%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.267"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -9.81%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.28"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -8.19%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "26.39"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -7.72%  "

$ws.Range("E48").Value = "  -0.03%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.32"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -5.52%  "

$ws.Range("E50").Value = "  -4.43%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "118.14"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
